# Generate Report for Handback
# Reorders the localization-status rows by file name and marks
# 56279519-...md and cbb81f83-...md as handed back (in sync with en-US),
# adding the "Latest Target File" / "Latest Handback File" / handback
# timestamp detail on the per-language sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet 1: "Overview"
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

# Clear existing hyperlinks on the sheet (any range works - the host
# clears the whole sheet's hyperlink collection) so we can rebuild them
# in the new row order without stale duplicates.
$ws1.Range("A1").Hyperlinks.Delete()

# Row 2: 56279519-...
$ws1.Cells.Item(2,2).Value = "Handed back: in sync with en-US"
$ws1.Cells.Item(2,3).Value = "Handed back: in sync with en-US"
$ws1.Cells.Item(2,4).Value = "2016-03-21 12:16:11"

# Row 3: cbb81f83-...
$ws1.Cells.Item(3,2).Value = "Handed back: in sync with en-US"
$ws1.Cells.Item(3,3).Value = "Handed back: in sync with en-US"
$ws1.Cells.Item(3,4).Value = "2016-03-21 12:16:11"

# Row 4: b6dd9f3c-...
$ws1.Cells.Item(4,2).Value = "In Translation"
$ws1.Cells.Item(4,3).Value = "In Translation"
$ws1.Cells.Item(4,4).Value = "2016-03-21 12:14:49"

# Row 5: ecee2a21-... (unchanged content, kept for clarity)
$ws1.Cells.Item(5,2).Value = "Ready for handoff"
$ws1.Cells.Item(5,3).Value = "Ready for handoff"
$ws1.Cells.Item(5,4).Value = "2016-03-21 12:16:11"

$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/bddcc38e0d431ec4ab74ba74f604a8b6f2bf548d/e2e/56279519-6b4a-4099-8473-fe409ec83634.md", "", "", "56279519-6b4a-4099-8473-fe409ec83634.md")
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/bddcc38e0d431ec4ab74ba74f604a8b6f2bf548d/e2e/cbb81f83-0ec0-4f3a-9148-7152febd4912.md", "", "", "cbb81f83-0ec0-4f3a-9148-7152febd4912.md")
$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/6e8bcc3727e6254bd897651116d6872a6f0d293e/e2e/b6dd9f3c-a03a-4825-bb52-8040edf1deca.md", "", "", "b6dd9f3c-a03a-4825-bb52-8040edf1deca.md")
$ws1.Hyperlinks.Add($ws1.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/d839d5c3ab862d5372a794ab32cf4fc18e5a2b15/e2e/ecee2a21-1227-4d8c-b52d-c7098e3d2b03.md", "", "", "ecee2a21-1227-4d8c-b52d-c7098e3d2b03.md")

# ---------------------------------------------------------------
# Sheet 2: "zh-cn"
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("A1").Hyperlinks.Delete()

# Row 2: 56279519-...
$ws2.Cells.Item(2,1).Value = "56279519-6b4a-4099-8473-fe409ec83634.md"
$ws2.Cells.Item(2,2).Value = ".md"
$ws2.Cells.Item(2,3).Value = "Handed back: in sync with en-US"
$ws2.Cells.Item(2,4).Value = "56279519-6b4a-4099-8473-fe409ec83634.4a4b55ef18b7c3ac50edf3ab054a1fb04df85497.zh-cn.xlf"
$ws2.Cells.Item(2,5).Value = "2016-03-21 12:16:05"
$ws2.Cells.Item(2,6).Value = "56279519-6b4a-4099-8473-fe409ec83634.md"
$ws2.Cells.Item(2,7).Value = "56279519-6b4a-4099-8473-fe409ec83634.4a4b55ef18b7c3ac50edf3ab054a1fb04df85497.zh-cn.xlf"
$ws2.Cells.Item(2,8).Value = "2016-03-21 12:16:38"
$ws2.Cells.Item(2,10).Value = "Include"

# Row 3: cbb81f83-...
$ws2.Cells.Item(3,1).Value = "cbb81f83-0ec0-4f3a-9148-7152febd4912.md"
$ws2.Cells.Item(3,2).Value = ".md"
$ws2.Cells.Item(3,3).Value = "Handed back: in sync with en-US"
$ws2.Cells.Item(3,4).Value = "cbb81f83-0ec0-4f3a-9148-7152febd4912.6fd00e16d2c98d80558cd8f4ce0f6a17dd00961a.zh-cn.xlf"
$ws2.Cells.Item(3,5).Value = "2016-03-21 12:16:05"
$ws2.Cells.Item(3,6).Value = "cbb81f83-0ec0-4f3a-9148-7152febd4912.md"
$ws2.Cells.Item(3,7).Value = "cbb81f83-0ec0-4f3a-9148-7152febd4912.6fd00e16d2c98d80558cd8f4ce0f6a17dd00961a.zh-cn.xlf"
$ws2.Cells.Item(3,8).Value = "2016-03-21 12:16:38"
$ws2.Cells.Item(3,10).Value = "Include"

# Row 4: b6dd9f3c-...
$ws2.Cells.Item(4,1).Value = "b6dd9f3c-a03a-4825-bb52-8040edf1deca.md"
$ws2.Cells.Item(4,2).Value = ".md"
$ws2.Cells.Item(4,3).Value = "In Translation"
$ws2.Cells.Item(4,4).Value = "b6dd9f3c-a03a-4825-bb52-8040edf1deca.3e94f4546f167d617e836689360f1a8f9b0b3a89.zh-cn.xlf"
$ws2.Cells.Item(4,5).Value = "2016-03-21 12:14:45"
$ws2.Cells.Item(4,8).Value = "0001-01-01 00:00:00"
$ws2.Cells.Item(4,10).Value = "Include"

# Row 5: ecee2a21-...
$ws2.Cells.Item(5,1).Value = "ecee2a21-1227-4d8c-b52d-c7098e3d2b03.md"
$ws2.Cells.Item(5,2).Value = ".md"
$ws2.Cells.Item(5,3).Value = "Ready for handoff"
$ws2.Cells.Item(5,4).Value = "ecee2a21-1227-4d8c-b52d-c7098e3d2b03.e9833127a94fd59149dc36d1e7351cd7e058a16b.zh-cn.xlf"
$ws2.Cells.Item(5,5).Value = "2016-03-21 12:16:05"
$ws2.Cells.Item(5,8).Value = "0001-01-01 00:00:00"
$ws2.Cells.Item(5,10).Value = "Include"

$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/bddcc38e0d431ec4ab74ba74f604a8b6f2bf548d/e2e/56279519-6b4a-4099-8473-fe409ec83634.md", "", "", "56279519-6b4a-4099-8473-fe409ec83634.md")
$ws2.Hyperlinks.Add($ws2.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c420fa1396a83eb40924fad12bf478309705596f/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/56279519-6b4a-4099-8473-fe409ec83634.4a4b55ef18b7c3ac50edf3ab054a1fb04df85497.zh-cn.xlf", "", "", "56279519-6b4a-4099-8473-fe409ec83634.4a4b55ef18b7c3ac50edf3ab054a1fb04df85497.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/bddcc38e0d431ec4ab74ba74f604a8b6f2bf548d/e2e/56279519-6b4a-4099-8473-fe409ec83634.md", "", "", "56279519-6b4a-4099-8473-fe409ec83634.md")
$ws2.Hyperlinks.Add($ws2.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c420fa1396a83eb40924fad12bf478309705596f/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/56279519-6b4a-4099-8473-fe409ec83634.4a4b55ef18b7c3ac50edf3ab054a1fb04df85497.zh-cn.xlf", "", "", "56279519-6b4a-4099-8473-fe409ec83634.4a4b55ef18b7c3ac50edf3ab054a1fb04df85497.zh-cn.xlf")

$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/bddcc38e0d431ec4ab74ba74f604a8b6f2bf548d/e2e/cbb81f83-0ec0-4f3a-9148-7152febd4912.md", "", "", "cbb81f83-0ec0-4f3a-9148-7152febd4912.md")
$ws2.Hyperlinks.Add($ws2.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c420fa1396a83eb40924fad12bf478309705596f/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/cbb81f83-0ec0-4f3a-9148-7152febd4912.6fd00e16d2c98d80558cd8f4ce0f6a17dd00961a.zh-cn.xlf", "", "", "cbb81f83-0ec0-4f3a-9148-7152febd4912.6fd00e16d2c98d80558cd8f4ce0f6a17dd00961a.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/bddcc38e0d431ec4ab74ba74f604a8b6f2bf548d/e2e/cbb81f83-0ec0-4f3a-9148-7152febd4912.md", "", "", "cbb81f83-0ec0-4f3a-9148-7152febd4912.md")
$ws2.Hyperlinks.Add($ws2.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c420fa1396a83eb40924fad12bf478309705596f/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/cbb81f83-0ec0-4f3a-9148-7152febd4912.6fd00e16d2c98d80558cd8f4ce0f6a17dd00961a.zh-cn.xlf", "", "", "cbb81f83-0ec0-4f3a-9148-7152febd4912.6fd00e16d2c98d80558cd8f4ce0f6a17dd00961a.zh-cn.xlf")

$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/6e8bcc3727e6254bd897651116d6872a6f0d293e/e2e/b6dd9f3c-a03a-4825-bb52-8040edf1deca.md", "", "", "b6dd9f3c-a03a-4825-bb52-8040edf1deca.md")
$ws2.Hyperlinks.Add($ws2.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/de2fcdc429e205cd22b651da168484252c9435e7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/b6dd9f3c-a03a-4825-bb52-8040edf1deca.3e94f4546f167d617e836689360f1a8f9b0b3a89.zh-cn.xlf", "", "", "b6dd9f3c-a03a-4825-bb52-8040edf1deca.3e94f4546f167d617e836689360f1a8f9b0b3a89.zh-cn.xlf")

$ws2.Hyperlinks.Add($ws2.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/d839d5c3ab862d5372a794ab32cf4fc18e5a2b15/e2e/ecee2a21-1227-4d8c-b52d-c7098e3d2b03.md", "", "", "ecee2a21-1227-4d8c-b52d-c7098e3d2b03.md")
$ws2.Hyperlinks.Add($ws2.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c420fa1396a83eb40924fad12bf478309705596f/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ecee2a21-1227-4d8c-b52d-c7098e3d2b03.e9833127a94fd59149dc36d1e7351cd7e058a16b.zh-cn.xlf", "", "", "ecee2a21-1227-4d8c-b52d-c7098e3d2b03.e9833127a94fd59149dc36d1e7351cd7e058a16b.zh-cn.xlf")

# ---------------------------------------------------------------
# Sheet 3: "de-de"
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("A1").Hyperlinks.Delete()

# Row 2: 56279519-...
$ws3.Cells.Item(2,1).Value = "56279519-6b4a-4099-8473-fe409ec83634.md"
$ws3.Cells.Item(2,2).Value = ".md"
$ws3.Cells.Item(2,3).Value = "Handed back: in sync with en-US"
$ws3.Cells.Item(2,4).Value = "56279519-6b4a-4099-8473-fe409ec83634.4a4b55ef18b7c3ac50edf3ab054a1fb04df85497.de-de.xlf"
$ws3.Cells.Item(2,5).Value = "2016-03-21 12:16:11"
$ws3.Cells.Item(2,6).Value = "56279519-6b4a-4099-8473-fe409ec83634.md"
$ws3.Cells.Item(2,7).Value = "56279519-6b4a-4099-8473-fe409ec83634.4a4b55ef18b7c3ac50edf3ab054a1fb04df85497.de-de.xlf"
$ws3.Cells.Item(2,8).Value = "2016-03-21 12:16:43"
$ws3.Cells.Item(2,10).Value = "Include"

# Row 3: cbb81f83-...
$ws3.Cells.Item(3,1).Value = "cbb81f83-0ec0-4f3a-9148-7152febd4912.md"
$ws3.Cells.Item(3,2).Value = ".md"
$ws3.Cells.Item(3,3).Value = "Handed back: in sync with en-US"
$ws3.Cells.Item(3,4).Value = "cbb81f83-0ec0-4f3a-9148-7152febd4912.6fd00e16d2c98d80558cd8f4ce0f6a17dd00961a.de-de.xlf"
$ws3.Cells.Item(3,5).Value = "2016-03-21 12:16:11"
$ws3.Cells.Item(3,6).Value = "cbb81f83-0ec0-4f3a-9148-7152febd4912.md"
$ws3.Cells.Item(3,7).Value = "cbb81f83-0ec0-4f3a-9148-7152febd4912.6fd00e16d2c98d80558cd8f4ce0f6a17dd00961a.de-de.xlf"
$ws3.Cells.Item(3,8).Value = "2016-03-21 12:16:43"
$ws3.Cells.Item(3,10).Value = "Include"

# Row 4: b6dd9f3c-...
$ws3.Cells.Item(4,1).Value = "b6dd9f3c-a03a-4825-bb52-8040edf1deca.md"
$ws3.Cells.Item(4,2).Value = ".md"
$ws3.Cells.Item(4,3).Value = "In Translation"
$ws3.Cells.Item(4,4).Value = "b6dd9f3c-a03a-4825-bb52-8040edf1deca.3e94f4546f167d617e836689360f1a8f9b0b3a89.de-de.xlf"
$ws3.Cells.Item(4,5).Value = "2016-03-21 12:14:49"
$ws3.Cells.Item(4,8).Value = "0001-01-01 00:00:00"
$ws3.Cells.Item(4,10).Value = "Include"

# Row 5: ecee2a21-...
$ws3.Cells.Item(5,1).Value = "ecee2a21-1227-4d8c-b52d-c7098e3d2b03.md"
$ws3.Cells.Item(5,2).Value = ".md"
$ws3.Cells.Item(5,3).Value = "Ready for handoff"
$ws3.Cells.Item(5,4).Value = "ecee2a21-1227-4d8c-b52d-c7098e3d2b03.e9833127a94fd59149dc36d1e7351cd7e058a16b.de-de.xlf"
$ws3.Cells.Item(5,5).Value = "2016-03-21 12:16:11"
$ws3.Cells.Item(5,8).Value = "0001-01-01 00:00:00"
$ws3.Cells.Item(5,10).Value = "Include"

$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/bddcc38e0d431ec4ab74ba74f604a8b6f2bf548d/e2e/56279519-6b4a-4099-8473-fe409ec83634.md", "", "", "56279519-6b4a-4099-8473-fe409ec83634.md")
$ws3.Hyperlinks.Add($ws3.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/bc9d1d4a91585f57215c9f2990de802661934012/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/56279519-6b4a-4099-8473-fe409ec83634.4a4b55ef18b7c3ac50edf3ab054a1fb04df85497.de-de.xlf", "", "", "56279519-6b4a-4099-8473-fe409ec83634.4a4b55ef18b7c3ac50edf3ab054a1fb04df85497.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/bddcc38e0d431ec4ab74ba74f604a8b6f2bf548d/e2e/56279519-6b4a-4099-8473-fe409ec83634.md", "", "", "56279519-6b4a-4099-8473-fe409ec83634.md")
$ws3.Hyperlinks.Add($ws3.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/bc9d1d4a91585f57215c9f2990de802661934012/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/56279519-6b4a-4099-8473-fe409ec83634.4a4b55ef18b7c3ac50edf3ab054a1fb04df85497.de-de.xlf", "", "", "56279519-6b4a-4099-8473-fe409ec83634.4a4b55ef18b7c3ac50edf3ab054a1fb04df85497.de-de.xlf")

$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/bddcc38e0d431ec4ab74ba74f604a8b6f2bf548d/e2e/cbb81f83-0ec0-4f3a-9148-7152febd4912.md", "", "", "cbb81f83-0ec0-4f3a-9148-7152febd4912.md")
$ws3.Hyperlinks.Add($ws3.Range("D3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/bc9d1d4a91585f57215c9f2990de802661934012/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/cbb81f83-0ec0-4f3a-9148-7152febd4912.6fd00e16d2c98d80558cd8f4ce0f6a17dd00961a.de-de.xlf", "", "", "cbb81f83-0ec0-4f3a-9148-7152febd4912.6fd00e16d2c98d80558cd8f4ce0f6a17dd00961a.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/bddcc38e0d431ec4ab74ba74f604a8b6f2bf548d/e2e/cbb81f83-0ec0-4f3a-9148-7152febd4912.md", "", "", "cbb81f83-0ec0-4f3a-9148-7152febd4912.md")
$ws3.Hyperlinks.Add($ws3.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/bc9d1d4a91585f57215c9f2990de802661934012/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/cbb81f83-0ec0-4f3a-9148-7152febd4912.6fd00e16d2c98d80558cd8f4ce0f6a17dd00961a.de-de.xlf", "", "", "cbb81f83-0ec0-4f3a-9148-7152febd4912.6fd00e16d2c98d80558cd8f4ce0f6a17dd00961a.de-de.xlf")

$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/6e8bcc3727e6254bd897651116d6872a6f0d293e/e2e/b6dd9f3c-a03a-4825-bb52-8040edf1deca.md", "", "", "b6dd9f3c-a03a-4825-bb52-8040edf1deca.md")
$ws3.Hyperlinks.Add($ws3.Range("D4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/392009882669f950c99c0e07dfb68707739dc44a/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/b6dd9f3c-a03a-4825-bb52-8040edf1deca.3e94f4546f167d617e836689360f1a8f9b0b3a89.de-de.xlf", "", "", "b6dd9f3c-a03a-4825-bb52-8040edf1deca.3e94f4546f167d617e836689360f1a8f9b0b3a89.de-de.xlf")

$ws3.Hyperlinks.Add($ws3.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/d839d5c3ab862d5372a794ab32cf4fc18e5a2b15/e2e/ecee2a21-1227-4d8c-b52d-c7098e3d2b03.md", "", "", "ecee2a21-1227-4d8c-b52d-c7098e3d2b03.md")
$ws3.Hyperlinks.Add($ws3.Range("D5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/bc9d1d4a91585f57215c9f2990de802661934012/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ecee2a21-1227-4d8c-b52d-c7098e3d2b03.e9833127a94fd59149dc36d1e7351cd7e058a16b.de-de.xlf", "", "", "ecee2a21-1227-4d8c-b52d-c7098e3d2b03.e9833127a94fd59149dc36d1e7351cd7e058a16b.de-de.xlf")
